$wb = $excel.ActiveWorkbook

# Add Sheet3 at the end of the workbook
$tmp = $wb.Worksheets.Add()
$tmp.Name = "Sheet3"
$tmp.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-fetch a live reference after the move
$ws3 = $wb.Worksheets.Item("Sheet3")

# Headers
$ws3.Range("A1").Value = "C1"
$ws3.Range("B1").Value = "C2"
$ws3.Range("C1").Value = "T"

# Data
$ws3.Range("A2").Value = 1
$ws3.Range("A3").Value = 2
$ws3.Range("A4").Value = 3
$ws3.Range("B2").Value = 4
$ws3.Range("B3").Value = 5
$ws3.Range("B4").Value = 6

# Formula
$ws3.Range("C4").Formula = "=SUM(A2:B4)"

# Make Sheet3 the active/selected sheet
$ws3.Select() | Out-Null
$ws3.Range("C4").Select() | Out-Null
